$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  +1.08%  '
$ws.Range("E8").Value = '  +2.29%  '
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("E11").Value = '  +5.37%  '
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("E21").Value = '  +2.56%  '
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +1.58%  '
$ws.Range("E27").Value = '  +0.84%  '
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("E29").Value = '  +0.98%  '
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("E31").Value = '  +4.54%  '
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("E33").Value = '  +4.18%  '
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("E37").Value = '  +3.45%  '
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("E40").Value = '  +2.54%  '
$ws.Range("E41").Value = '  +1.47%  '
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("E43").Value = '  +2.33%  '
$ws.Range("E44").Value = '  +4.19%  '
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("E47").Value = '  +3.35%  '
$ws.Range("E48").Value = '  +2.01%  '
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("E51").Value = '  +4.21%  '
